$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.440.53"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").Value = "1.918.35"
$ws.Range("E3").Value = "  +0.73%  "

$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.50%  "

$ws.Range("D5").Value = "'325.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.12%  "

$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("D7").Value = "'0.4823"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.48%  "

$ws.Range("D8").Value = "'0.4069"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.08218"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.73%  "

$ws.Range("E10").Value = "  +1.07%  "

$ws.Range("D11").Value = "'23.28"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.67%  "

$ws.Range("D12").Value = "1.942.65"
$ws.Range("E12").Value = "  +2.42%  "

$ws.Range("D13").Value = "'6.069"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.89%  "

$ws.Range("D14").Value = "'7.249"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.42%  "

$ws.Range("D15").Value = "'91.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.91%  "

$ws.Range("D16").Value = "'0.06875"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.84%  "

$ws.Range("D17").Value = "'1.009"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.50%  "

$ws.Range("E18").Value = "  +0.56%  "

$ws.Range("D19").Value = "'17.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("E20").Value = "  +0.57%  "

$ws.Range("D21").Value = "29.461.61"
$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("E22").Value = "  +2.25%  "

$ws.Range("D23").Value = "'11.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.44%  "

$ws.Range("D24").Value = "'2.188"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.03%  "

$ws.Range("D25").Value = "2.159.84"
$ws.Range("E25").Value = "  +1.60%  "

$ws.Range("D26").Value = "'6.668"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.32%  "

$ws.Range("D27").Value = "'156.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.79%  "

$ws.Range("D28").Value = "'20.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.62%  "

$ws.Range("D29").Value = "'2.118"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.90%  "

$ws.Range("D30").Value = "'120.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.14%  "

$ws.Range("D31").Value = "'1.015"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.96%  "

$ws.Range("D32").Value = "'0.09618"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.16%  "

$ws.Range("D33").Value = "'5.648"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.22%  "

$ws.Range("D34").Value = "'3.548"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.13%  "

$ws.Range("D35").Value = "'1.375"
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").Value = "'0.02286"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.47%  "

$ws.Range("E37").Value = "  +0.39%  "

$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'10.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.84%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.066"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.94%  "

$ws.Range("E41").Value = "  +1.20%  "

$ws.Range("D42").Value = "'0.1848"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("E43").Value = "  +0.16%  "

$ws.Range("D44").Value = "'2.405"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.20%  "

$ws.Range("D45").Value = "'0.07601"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.69%  "

$ws.Range("E46").Value = "  +1.11%  "

$ws.Range("D47").Value = "'0.5597"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.22%  "

$ws.Range("D48").Value = "'1.956"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.57%  "

$ws.Range("D49").Value = "'118.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.97%  "

$ws.Range("D50").Value = "'2.427"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.69%  "

$ws.Range("D51").Value = "'72.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.09%  "
